# Append a new search-result row ("Noun" method run) to the tracking sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 7
$ws.Range("A$row").Value = 42602.582546296297
$ws.Range("B$row").Value = "Noun"
$ws.Range("C$row").Value = 7724
$ws.Range("D$row").Value = 6151
$ws.Range("E$row").Value = 1145
$ws.Range("F$row").Value = 133
$ws.Range("G$row").Value = 61
$ws.Range("H$row").Value = 67
$ws.Range("I$row").Value = 31
$ws.Range("J$row").Value = 3
$ws.Range("K$row").Value = 3
$ws.Range("L$row").Value = 49
$ws.Range("M$row").Value = 49
